$wb = $excel.ActiveWorkbook

# 1. Update the status text "Ready for handoff" -> "In Translation" on every
#    sheet that references it (Overview!E2:F2, zh-cn!C2, de-de!C2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# 2. Narrow the "zh-cn"/"de-de" columns on Overview (E, F) and the "Status"
#    column on the language sheets (C) from 17.2159881591797 to
#    13.4101845877511 (characters). Excel's ColumnWidth setter snaps to
#    whole-pixel increments, so feed it the pre-offset "characters" value
#    that lands on the pixel closest to the target stored width.
$targetWidth = 12.5
$wsOverview.Range("E1").ColumnWidth = $targetWidth
$wsOverview.Range("F1").ColumnWidth = $targetWidth

$wsZhCn.Range("C1").ColumnWidth = $targetWidth
$wsDeDe.Range("C1").ColumnWidth = $targetWidth
